$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Failed"
$ws.Range("H2").Value = "17/04/2021"
$ws.Range("H3").Value = "17/04/2021"
